$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1725.2941
$ws.Range("I2").Value = 283.75
$ws.Range("J2").Value = 3006.6667
$ws.Range("K2").Value = 283.75
$ws.Range("L2").Value = 3006.6667
$ws.Range("M2").Value = -170.75
$ws.Range("N2").Value = -3232.6667
# Row 32
$ws.Range("H32").Value = 4011.25
$ws.Range("I32").Value = 1896.6666
$ws.Range("J32").Value = 5280
$ws.Range("K32").Value = 1896.6666
$ws.Range("L32").Value = 5280
$ws.Range("M32").Value = -1570.6666
$ws.Range("N32").Value = -5932
# Row 70
$ws.Range("H70").Value = 1557.1666
$ws.Range("I70").Value = 998
$ws.Range("K70").Value = 2994
$ws.Range("M70").Value = -2724
# Row 73
$ws.Range("H73").Value = 1557.1666
$ws.Range("I73").Value = 998
$ws.Range("K73").Value = 2994
$ws.Range("M73").Value = -2058
# Row 88
$ws.Range("H88").Value = 2374
$ws.Range("J88").Value = 2374
$ws.Range("L88").Value = 2374
$ws.Range("N88").Value = -3186
# Row 91
$ws.Range("H91").Value = 2374
$ws.Range("J91").Value = 2374
$ws.Range("L91").Value = 2374
$ws.Range("N91").Value = -5182
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
# Row 127
$ws.Range("H127").Value = 1678.4445
$ws.Range("I127").Value = 2184.8
$ws.Range("K127").Value = 6554.400000000001
$ws.Range("M127").Value = -1594.400000000001
# Row 138
$ws.Range("H138").Value = 4318.2666
$ws.Range("J138").Value = 4694.3477
$ws.Range("L138").Value = 14083.0431
$ws.Range("N138").Value = -24363.0431

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 725
$ws.Range("I74").Value = 725
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 725
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 149
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 725
$ws.Range("I77").Value = 725
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3625
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 743
$ws.Range("N77").ClearContents()
# Row 88
$ws.Range("H88").Value = 3132.5454
$ws.Range("J88").Value = 3224.25
$ws.Range("L88").Value = 3224.25
$ws.Range("N88").Value = -4036.25
# Row 91
$ws.Range("H91").Value = 3132.5454
$ws.Range("J91").Value = 3224.25
$ws.Range("L91").Value = 3224.25
$ws.Range("N91").Value = -6032.25
# Row 122
$ws.Range("H122").Value = 2246.875
$ws.Range("I122").Value = 1910.7142
$ws.Range("K122").Value = 5732.142599999999
$ws.Range("M122").Value = -3282.142599999999
# Row 132
$ws.Range("H132").Value = 2646
$ws.Range("I132").Value = 2646
$ws.Range("K132").Value = 7938
$ws.Range("M132").Value = -5408

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 23499.8
$ws.Range("J33").Value = 24375
$ws.Range("L33").Value = 24375
$ws.Range("N33").Value = -25047
# Row 105
$ws.Range("H105").Value = 2669.4285
$ws.Range("I105").Value = 2747
$ws.Range("J105").Value = 2566
$ws.Range("K105").Value = 2747
$ws.Range("L105").Value = 2566
$ws.Range("M105").Value = -1000
$ws.Range("N105").Value = -6060
# Row 134
$ws.Range("H134").Value = 2167
$ws.Range("I134").Value = 2218.111
$ws.Range("J134").Value = 1707
$ws.Range("K134").Value = 6654.333
$ws.Range("L134").Value = 5121
$ws.Range("M134").Value = -4119.333
$ws.Range("N134").Value = -10191

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2014.125
$ws.Range("I58").Value = 1998.2
$ws.Range("K58").Value = 1998.2
$ws.Range("M58").Value = -1795.2
# Row 107
$ws.Range("H107").Value = 1179.7778
$ws.Range("I107").Value = 1161.5714
$ws.Range("K107").Value = 1161.5714
$ws.Range("M107").Value = 758.4286
# Row 132
$ws.Range("H132").Value = 2441.4075
$ws.Range("I132").Value = 2573.238
$ws.Range("J132").Value = 1980
$ws.Range("K132").Value = 7719.714
$ws.Range("L132").Value = 5940
$ws.Range("M132").Value = -5189.714
$ws.Range("N132").Value = -11000
# Row 134
$ws.Range("H134").Value = 1868.4375
$ws.Range("I134").Value = 1828.2142
$ws.Range("K134").Value = 5484.642599999999
$ws.Range("M134").Value = -2949.642599999999
# Row 136
$ws.Range("H136").Value = 2014.125
$ws.Range("I136").Value = 1998.2
$ws.Range("K136").Value = 5994.6
$ws.Range("M136").Value = -3444.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 199.5
$ws.Range("I8").Value = 199.5
$ws.Range("K8").Value = 598.5
$ws.Range("M8").Value = -459.5
# Row 12
$ws.Range("H12").Value = 290.25
$ws.Range("J12").Value = 260.5
$ws.Range("L12").Value = 781.5
$ws.Range("N12").Value = -1127.5
# Row 109
$ws.Range("H109").Value = 1725.2858
$ws.Range("J109").Value = 4900
$ws.Range("L109").Value = 14700
$ws.Range("N109").Value = -16780
# Row 121
$ws.Range("H121").Value = 8519.6
$ws.Range("J121").Value = 3237.1765
$ws.Range("L121").Value = 9711.529500000001
$ws.Range("N121").Value = -12331.5295
# Row 141
$ws.Range("H141").Value = 8792.857
$ws.Range("I141").Value = 9425
$ws.Range("K141").Value = 28275
$ws.Range("M141").Value = -23095

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
# Row 132
$ws.Range("H132").Value = 2127.0667
$ws.Range("I132").Value = 2127.0667
$ws.Range("K132").Value = 6381.2001
$ws.Range("M132").Value = -3851.2001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 779.1667
$ws.Range("I16").Value = 795
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 795
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -625
$ws.Range("N16").Value = -1040
# Row 20
$ws.Range("H20").Value = 8712.5
$ws.Range("I20").Value = 8712.5
$ws.Range("K20").Value = 8712.5
$ws.Range("M20").Value = -8486.5
# Row 40
$ws.Range("H40").Value = 3332.3333
$ws.Range("I40").Value = 2248.5
$ws.Range("K40").Value = 2248.5
$ws.Range("M40").Value = -2112.5
# Row 55
$ws.Range("H55").Value = 2600.5715
$ws.Range("I55").Value = 2425.5
$ws.Range("J55").Value = 2834
$ws.Range("K55").Value = 2425.5
$ws.Range("L55").Value = 2834
$ws.Range("M55").Value = -2252.5
$ws.Range("N55").Value = -3180
# Row 61
$ws.Range("H61").Value = 1711.4445
$ws.Range("I61").Value = 1711.4445
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1711.4445
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1509.4445
$ws.Range("N61").ClearContents()
# Row 68
$ws.Range("H68").Value = 2300
$ws.Range("I68").Value = 2300
$ws.Range("K68").Value = 2300
$ws.Range("M68").Value = -1551
# Row 71
$ws.Range("H71").Value = 2300
$ws.Range("I71").Value = 2300
$ws.Range("K71").Value = 11500
$ws.Range("M71").Value = -7756
# Row 113
$ws.Range("H113").Value = 1711.4445
$ws.Range("I113").Value = 1711.4445
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1711.4445
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 458.5554999999999
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 5638.722
$ws.Range("I132").Value = 6958.4546
$ws.Range("J132").Value = 3564.8572
$ws.Range("K132").Value = 20875.3638
$ws.Range("L132").Value = 10694.5716
$ws.Range("M132").Value = -18345.3638
$ws.Range("N132").Value = -15754.5716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
# Row 75
$ws.Range("H75").Value = 71908
$ws.Range("I75").Value = 71908
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 71908
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -70972
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 71908
$ws.Range("I78").Value = 71908
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 215724
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -211044
$ws.Range("N78").ClearContents()
# Row 81
$ws.Range("H81").Value = 2001320.2
$ws.Range("I81").Value = 1400
$ws.Range("J81").Value = 3334600.2
$ws.Range("K81").Value = 2800
$ws.Range("L81").Value = 6669200.4
$ws.Range("M81").Value = -1739
$ws.Range("N81").Value = -6671322.4
# Row 84
$ws.Range("H84").Value = 2001320.2
$ws.Range("I84").Value = 1400
$ws.Range("J84").Value = 3334600.2
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 33346002
$ws.Range("M84").Value = -8696
$ws.Range("N84").Value = -33356610
# Row 122
$ws.Range("H122").Value = 1858.6666
$ws.Range("I122").Value = 1858.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5575.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3125.9998
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 5214.069
$ws.Range("I126").Value = 4882.316
$ws.Range("J126").Value = 5844.4
$ws.Range("K126").Value = 14646.948
$ws.Range("L126").Value = 17533.2
$ws.Range("M126").Value = -12176.948
$ws.Range("N126").Value = -22473.2
# Row 132
$ws.Range("H132").Value = 2173.3333
$ws.Range("I132").Value = 1789.3572
$ws.Range("K132").Value = 5368.071599999999
$ws.Range("M132").Value = -2838.071599999999
